# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — copy the header formatting (bold/centered/bordered) from
# the neighboring "Unnamed: 28" header cell, then set the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-42 — every player row gets the team's season record.
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 66  # AD
    $ws.Cells.Item($row, 31).Value = 47  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
